# Auto-generated edit script: updates cryptos.xlsx cell values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.800.09'
$ws.Range('E2').Value = '  +1.10%  '
$ws.Range('D3').Value = '2.306.73'
$ws.Range('E3').Value = '  +0.89%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''498.90'
$ws.Range('E5').Value = '  +1.21%  '
$ws.Range('D6').Value = '''129.21'
$ws.Range('E6').Value = '  +0.98%  '
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '''0.532'
$ws.Range('E8').Value = '  +0.93%  '
$ws.Range('D9').Value = '2.306.09'
$ws.Range('E9').Value = '  +0.67%  '
$ws.Range('D10').Value = '''0.0956'
$ws.Range('E10').Value = '  +1.99%  '
$ws.Range('E11').Value = '  +2.35%  '
$ws.Range('E12').Value = '  +2.30%  '
$ws.Range('D13').Value = '''4.66'
$ws.Range('E13').Value = '  -1.98%  '
$ws.Range('D14').Value = '2.713.36'
$ws.Range('E14').Value = '  +0.77%  '
$ws.Range('E15').Value = '  +2.79%  '
$ws.Range('D16').Value = '54.741.56'
$ws.Range('E16').Value = '  +1.02%  '
$ws.Range('E17').Value = '  +0.70%  '
$ws.Range('D18').Value = '2.285.07'
$ws.Range('E18').Value = '  +0.70%  '
$ws.Range('D19').Value = '''10.11'
$ws.Range('E19').Value = '  +4.42%  '
$ws.Range('D20').Value = '''4.11'
$ws.Range('E20').Value = '  +2.67%  '
$ws.Range('D21').Value = '''308.10'
$ws.Range('E21').Value = '  +1.60%  '
$ws.Range('D22').Value = '''6.47'
$ws.Range('E22').Value = '  +4.73%  '
$ws.Range('D23').Value = '''0.999'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('E24').Value = '  -1.13%  '
$ws.Range('D25').Value = '''62.96'
$ws.Range('E25').Value = '  -1.77%  '
$ws.Range('D26').Value = '''0.999'
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').Value = '''0.152'
$ws.Range('E27').Value = '  +5.87%  '
$ws.Range('B28').Value = 'Polygon'
$ws.Range('C28').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D28').Value = '''0.375'
$ws.Range('E28').Value = '  +2.34%  '
$ws.Range('D29').Value = '2.391.18'
$ws.Range('D30').Value = '''7.19'
$ws.Range('E30').Value = '  +0.87%  '
$ws.Range('D31').Value = '''170.11'
$ws.Range('E31').Value = '  +0.37%  '
$ws.Range('D32').Value = '0.0₃0702'
$ws.Range('E32').Value = '  +0.19%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('E34').Value = '  +3.05%  '
$ws.Range('D35').Value = '''0.999'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E36').Value = '  +1.18%  '
$ws.Range('D37').Value = '''0.999'
$ws.Range('E37').Value = '  +0.16%  '
$ws.Range('D38').Value = '''17.72'
$ws.Range('E38').Value = '  +0.60%  '
$ws.Range('E39').Value = '  +3.13%  '
$ws.Range('D40').Value = '''0.870'
$ws.Range('E40').Value = '  +1.97%  '
$ws.Range('E41').Value = '  +1.92%  '
$ws.Range('D42').Value = '''35.54'
$ws.Range('E42').Value = '  -0.63%  '
$ws.Range('E43').Value = '  +2.54%  '
$ws.Range('E44').Value = '  +2.35%  '
$ws.Range('D45').Value = '''3.37'
$ws.Range('E45').Value = '  +1.16%  '
$ws.Range('D46').Value = '''128.58'
$ws.Range('E46').Value = '  +4.14%  '
$ws.Range('D47').Value = '''4.92'
$ws.Range('E47').Value = '  +2.49%  '
$ws.Range('D48').Value = '''0.0895'
$ws.Range('E48').Value = '  +1.40%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').Value = '''245.65'
$ws.Range('E49').Value = '  +2.11%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '''0.552'
$ws.Range('E50').Value = '  +1.38%  '
$ws.Range('D51').Value = '''0.0487'
$ws.Range('E51').Value = '  +2.62%  '
